# Clean up the "Stock Trades" sheet:
#  - strip stray leading/trailing/internal whitespace baked into the
#    string cells (ticker codes, company names, sector names, headers)
#  - turn the free-text "Trade Date" column into real Excel date values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "StockCode"
$ws.Cells.Item(1, 2).Value = "CompanyName"
$ws.Cells.Item(1, 3).Value = "Sector"
$ws.Cells.Item(1, 4).Value = "Open"
$ws.Cells.Item(1, 5).Value = "Close"
$ws.Cells.Item(1, 6).Value = "Volume"
$ws.Cells.Item(1, 7).Value = "TradeDate"
$ws.Cells.Item(1, 8).Value = "MarketCap"

# Data rows: StockCode, CompanyName, Sector, Open, Close, Volume, TradeDate, MarketCap
$rows = @(
    @("BHP",   "BHPGroup",          "Materials",   4.2,               45.9,    2500000, "3/2/2024", 100),
    @("CBA",   "CommonwealthBank",  "Financials",  102.5,             100,     1850000, "3/1/2024", 175.2),
    @("WBC",   "WestpacBanking",    "Financials",  23,                4.1,     12000,   "3/1/2024", 85.3),
    @("CSL",   "CSLLimited",        "Financials",  291.4,             292.8,   455000,  "3/5/2024", 141),
    @("RIO",   "RioTintoLimited",   "Materials",   127.8,             128.15,  680000,  "3/1/2024", 47.7),
    @("ANZ",   "ANZBankingGroup",   "Healthcare",  28.9,              2.15,    950000,  "3/1/2024", 2.1),
    @("Achal", "ANZBankingGroup",   "Healthcare",  28.9,              2.15,    950000,  "3/1/2024", 2.1)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 7).NumberFormat = "m/d/yyyy"
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$ws.Range("C10").Select()
